$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; this shifts the former row 3 (and every
# row below it, including the former row 106) down by one, to rows 4..107.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the descriptive data copied from the row
# that is now directly below it (row 4, the former row 3), since this is
# the same market/product record, just a different reporting date.
$ws.Range("A3").Value2 = $ws.Range("A4").Value2
$ws.Range("B3").Value2 = $ws.Range("B4").Value2
$ws.Range("C3").Value2 = $ws.Range("C4").Value2
$ws.Range("E3").Value2 = $ws.Range("E4").Value2
$ws.Range("F3").Value2 = $ws.Range("F4").Value2
$ws.Range("G3").Value2 = $ws.Range("G4").Value2
$ws.Range("H3").Value2 = $ws.Range("H4").Value2
$ws.Range("I3").Value2 = $ws.Range("I4").Value2
$ws.Range("J3").Value2 = $ws.Range("J4").Value2
$ws.Range("K3").Value2 = $ws.Range("K4").Value2
$ws.Range("L3").Value2 = $ws.Range("L4").Value2
$ws.Range("Q3").Value2 = $ws.Range("Q4").Value2
$ws.Range("R3").Value2 = $ws.Range("R4").Value2
$ws.Range("T3").Value2 = $ws.Range("T4").Value2

# Now set the new values that are specific to this newly reported record.
$ws.Range("D3").Value2 = 44599
$ws.Range("M3").Value2 = 240
$ws.Range("N3").Value2 = 7000
$ws.Range("O3").Value2 = 7000
$ws.Range("P3").Value2 = 7000
$ws.Range("S3").Value2 = 1750
